# "Made use of merged column reading"
#
# The SetEffects column used to hold one big comma/bracket separated blob
# like:
#   ,
#   [StatisticBoost;DamageModifier_PhysicalDamage=0.25;],
#   ,
#   [StatisticBoost,Triggerable;DamageModifier_ChargedAttack=0.25;],
#
# This rewrites the sheet so the 2pc / 4pc set bonuses live in their own
# cells (D = 2pc bonus, F = 4pc bonus) under a single merged "SetEffects"
# header spanning C:G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlHAlignCenter = -4108
$xlHAlignLeft   = -4131
$xlVAlignCenter = -4108

# ---------------------------------------------------------------------
# Columns: C used to be one very wide column (135.5) holding the whole
# blob; now C:G are each a more modest 76.375 wide.
# ---------------------------------------------------------------------
$ws.Range("C1:G1").ColumnWidth = 76.375

# ---------------------------------------------------------------------
# Row 2 header: "SetEffects" becomes a merged header cell C2:G2, and the
# three header cells are centered / left-aligned instead of the default.
# ---------------------------------------------------------------------
$ws.Range("A2:B2").HorizontalAlignment = $xlHAlignLeft
$ws.Range("A2:B2").VerticalAlignment = $xlVAlignCenter

$ws.Range("C2:G2").HorizontalAlignment = $xlHAlignLeft
$ws.Range("C2:G2").VerticalAlignment = $xlVAlignCenter

$ws.Range("C2:G2").Merge()

# ---------------------------------------------------------------------
# Row 3 (Pale Flame): split the old combined blob into 2pc bonus (D3)
# and 4pc bonus (F3); clear the old C3 text but keep its wrap style.
# ---------------------------------------------------------------------
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "StatisticBoost;DamageModifier_PhysicalDamage=0.25;"
$ws.Range("F3").Value = "StatisticBoost,Stackable;Attack_Percentage=0.09,MaxNumStacks=2,InitialNumStacks=0;,`nStatisticBoost,Triggerable;DamageModifier_PhysicalDamage=0.25;"
$ws.Range("F3").WrapText = $true
$ws.Rows(3).RowHeight = 29.25

# ---------------------------------------------------------------------
# Row 4 (Bloodstained Chivalry): same split for 2pc / 4pc bonus.
# ---------------------------------------------------------------------
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "StatisticBoost;DamageModifier_PhysicalDamage=0.25;"
$ws.Range("F4").Value = "StatisticBoost,Triggerable;DamageModifier_ChargedAttack=0.25;"
$ws.Rows(4).AutoFit()

# ---------------------------------------------------------------------
# Selection ends up parked on D11, matching the author's last click.
# ---------------------------------------------------------------------
$ws.Range("D11").Select()
